$wb = $excel.ActiveWorkbook

# Add a new worksheet named "BVT" after the existing "Meganav" sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "BVT"

# Header row (row 1)
$ws.Range("A1").Value = "Test"
$ws.Range("B1").Value = "Description"
$ws.Range("C1").Value = "Meganav1"
$ws.Range("D1").Value = "Expected Header"
$ws.Range("E1").Value = "Price Filter"
$ws.Range("F1").Value = "Color Filter"
$ws.Range("G1").Value = "Expected breadcrumb"
$ws.Range("H1").Value = "Test Rail url"

# Data row (row 2)
$ws.Range("A2").Value = "bvt-01"
$ws.Range("B2").Value = "Guest user should be able to add products to cart and checkout"
$ws.Range("C2").Value = "/c-28-serveware-flatware"
$ws.Range("D2").Value = "Serveware & Flatware"
$ws.Range("E2").Value = "`$60-`$79"
$ws.Range("F2").Value = "Gold"
$ws.Range("G2").Value = "Home  Tableware  Serveware & Flatware `$60 - `$79x  Clear All"

$testRailUrl = "https://surlatable.testrail.net/index.php?/cases/view/12080&group_by=cases:section_id&group_order=asc&display_deleted_cases=0&group_id=1961"
$ws.Range("H2").Value = $testRailUrl
$ws.Hyperlinks.Add($ws.Range("H2"), $testRailUrl) | Out-Null

# Column widths (characters), tuned so the engine's stored (rounded) width
# lands as close as possible to the authored workbook's best-fit widths
# (11.7109375, 57.7109375, 31.42578125, 30.85546875, 30.7109375, 14.140625,
# 55.5703125, 138 characters respectively).
$ws.Columns.Item(1).ColumnWidth = 10.833333333333334
$ws.Columns.Item(2).ColumnWidth = 56.833333333333336
$ws.Columns.Item(3).ColumnWidth = 30.666666666666668
$ws.Columns.Item(4).ColumnWidth = 30.0
$ws.Columns.Item(5).ColumnWidth = 29.833333333333332
$ws.Columns.Item(6).ColumnWidth = 13.333333333333334
$ws.Columns.Item(7).ColumnWidth = 54.666666666666664
$ws.Columns.Item(8).ColumnWidth = 137.16666666666666

# Selection / active cell on the new sheet, and make it the active tab.
$ws.Range("D1").Select() | Out-Null
$ws.Activate()

Write-Output "Added BVT sheet with test record data"
